# Updates the cryptos list (Price / Volume(1h) columns, and two rank swaps)
# to match the latest scrape, matching the commit's unified diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells hold text that can look numeric (e.g. "1.00", "3.542.35")
# or percentage strings with padding spaces. Force text format before writing
# so Excel does not auto-coerce to a number/date, then clear the format back
# to the default (General) so no stray cell styling is introduced.
$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "65.117.14"
$rng.ClearFormats()
$rng = $ws.Range("E2")
$rng.NumberFormat = "@"
$rng.Value = "  +0.51%  "
$rng.ClearFormats()
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "3.568.00"
$rng.ClearFormats()
$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = "  +3.98%  "
$rng.ClearFormats()
$rng = $ws.Range("D4")
$rng.NumberFormat = "@"
$rng.Value = "0.999"
$rng.ClearFormats()
$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$rng.Value = "  -0.08%  "
$rng.ClearFormats()
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "597.80"
$rng.ClearFormats()
$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$rng.Value = "  +2.85%  "
$rng.ClearFormats()
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = "138.18"
$rng.ClearFormats()
$rng = $ws.Range("E6")
$rng.NumberFormat = "@"
$rng.Value = "  +2.76%  "
$rng.ClearFormats()
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "3.567.81"
$rng.ClearFormats()
$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$rng.Value = "  +3.99%  "
$rng.ClearFormats()
$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$rng.Value = "  +0.12%  "
$rng.ClearFormats()
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "0.495"
$rng.ClearFormats()
$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = "  +2.47%  "
$rng.ClearFormats()
$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = "  +2.68%  "
$rng.ClearFormats()
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "6.98"
$rng.ClearFormats()
$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$rng.Value = "  -1.43%  "
$rng.ClearFormats()
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "0.387"
$rng.ClearFormats()
$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$rng.Value = "  +3.55%  "
$rng.ClearFormats()
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "4.163.86"
$rng.ClearFormats()
$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$rng.Value = "  +3.76%  "
$rng.ClearFormats()
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "0.0000183"
$rng.ClearFormats()
$rng = $ws.Range("E14")
$rng.NumberFormat = "@"
$rng.Value = "  +2.59%  "
$rng.ClearFormats()
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "27.43"
$rng.ClearFormats()
$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$rng.Value = "  +4.74%  "
$rng.ClearFormats()
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "3.555.79"
$rng.ClearFormats()
$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$rng.Value = "  +3.40%  "
$rng.ClearFormats()
$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$rng.Value = "  +1.11%  "
$rng.ClearFormats()
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "64.955.61"
$rng.ClearFormats()
$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$rng.Value = "  +0.29%  "
$rng.ClearFormats()
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "10.21"
$rng.ClearFormats()
$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = "  +7.74%  "
$rng.ClearFormats()
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "5.88"
$rng.ClearFormats()
$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$rng.Value = "  +2.34%  "
$rng.ClearFormats()
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "14.34"
$rng.ClearFormats()
$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = "  +6.14%  "
$rng.ClearFormats()
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "391.25"
$rng.ClearFormats()
$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$rng.Value = "  +2.72%  "
$rng.ClearFormats()
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "0.579"
$rng.ClearFormats()
$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = "  +6.77%  "
$rng.ClearFormats()
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "3.698.62"
$rng.ClearFormats()
$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$rng.Value = "  +3.72%  "
$rng.ClearFormats()
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "74.13"
$rng.ClearFormats()
$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$rng.Value = "  +3.25%  "
$rng.ClearFormats()
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "1.00"
$rng.ClearFormats()
$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$rng.Value = "  +0.08%  "
$rng.ClearFormats()
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "0.0000116"
$rng.ClearFormats()
$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$rng.Value = "  +10.64%  "
$rng.ClearFormats()
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "7.82"
$rng.ClearFormats()
$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = "  +9.21%  "
$rng.ClearFormats()
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "0.999"
$rng.ClearFormats()
$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = "  -0.19%  "
$rng.ClearFormats()
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "2.29"
$rng.ClearFormats()
$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$rng.Value = "  +5.11%  "
$rng.ClearFormats()
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "8.33"
$rng.ClearFormats()
$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$rng.Value = "  +4.53%  "
$rng.ClearFormats()
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "3.569.60"
$rng.ClearFormats()
$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$rng.Value = "  +3.54%  "
$rng.ClearFormats()
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "1.45"
$rng.ClearFormats()
$rng = $ws.Range("E33")
$rng.NumberFormat = "@"
$rng.Value = "  +23.21%  "
$rng.ClearFormats()
$rng = $ws.Range("B34")
$rng.NumberFormat = "@"
$rng.Value = "EthereumClassic"
$rng.ClearFormats()
$rng = $ws.Range("C34")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$rng.ClearFormats()
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "24.00"
$rng.ClearFormats()
$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$rng.Value = "  +4.64%  "
$rng.ClearFormats()
$rng = $ws.Range("B35")
$rng.NumberFormat = "@"
$rng.Value = "USDe"
$rng.ClearFormats()
$rng = $ws.Range("C35")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$rng.ClearFormats()
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "1.00"
$rng.ClearFormats()
$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$rng.Value = "  -0.02%  "
$rng.ClearFormats()
$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$rng.Value = "  +1.81%  "
$rng.ClearFormats()
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "169.97"
$rng.ClearFormats()
$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$rng.Value = "  +0.31%  "
$rng.ClearFormats()
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "1.56"
$rng.ClearFormats()
$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$rng.Value = "  +8.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = "6.91"
$rng.ClearFormats()
$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$rng.Value = "  +3.19%  "
$rng.ClearFormats()
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "5.05"
$rng.ClearFormats()
$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$rng.Value = "  +10.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "0.0812"
$rng.ClearFormats()
$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$rng.Value = "  +7.14%  "
$rng.ClearFormats()
$rng = $ws.Range("B42")
$rng.NumberFormat = "@"
$rng.Value = "EnergySwap"
$rng.ClearFormats()
$rng = $ws.Range("C42")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$rng.ClearFormats()
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "27.02"
$rng.ClearFormats()
$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = "  +20.77%  "
$rng.ClearFormats()
$rng = $ws.Range("B43")
$rng.NumberFormat = "@"
$rng.Value = "Mantle"
$rng.ClearFormats()
$rng = $ws.Range("C43")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$rng.ClearFormats()
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "0.827"
$rng.ClearFormats()
$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$rng.Value = "  +3.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = "43.01"
$rng.ClearFormats()
$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$rng.Value = "  +2.58%  "
$rng.ClearFormats()
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "0.999"
$rng.ClearFormats()
$rng = $ws.Range("E45")
$rng.NumberFormat = "@"
$rng.Value = "  -0.20%  "
$rng.ClearFormats()
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "1.23"
$rng.ClearFormats()
$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$rng.Value = "  +9.98%  "
$rng.ClearFormats()
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "4.47"
$rng.ClearFormats()
$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = "  +4.48%  "
$rng.ClearFormats()
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "1.68"
$rng.ClearFormats()
$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$rng.Value = "  +4.44%  "
$rng.ClearFormats()
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "2.473.35"
$rng.ClearFormats()
$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$rng.Value = "  +12.89%  "
$rng.ClearFormats()
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = "6.92"
$rng.ClearFormats()
$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$rng.Value = "  +7.03%  "
$rng.ClearFormats()
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "302.87"
$rng.ClearFormats()
$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$rng.Value = "  +10.60%  "
$rng.ClearFormats()
